# Add 45C thermal-curve data (plate1, 2025-08-20) for all wells A01-A12,
# B01-B12 (sample) and C01-C06 (blank), appended after the existing data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$wells = @(
    "A01","A02","A03","A04","A05","A06","A07","A08","A09","A10","A11","A12",
    "B01","B02","B03","B04","B05","B06","B07","B08","B09","B10","B11","B12",
    "C01","C02","C03","C04","C05","C06"
)

$date = 20250820
$plate = "plate1"
$temperature = 45

$startRow = $ws.Cells.Item($ws.Rows.Count(), 1).End(-4162).Row() + 1

for ($i = 0; $i -lt $wells.Count; $i++) {
    $r = $startRow + $i
    $well = $wells[$i]
    if ($well.Substring(0, 1) -eq "C") {
        $type = "blank"
    } else {
        $type = "sample"
    }

    $ws.Cells.Item($r, 1).Value = $date
    $ws.Cells.Item($r, 2).Value = $plate
    $ws.Cells.Item($r, 3).Value = $temperature
    $ws.Cells.Item($r, 4).Value = $well
    $ws.Cells.Item($r, 5).Value = $type
}

$ws.Range("C157").Select()
